$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.223194389186745
$ws.Range("C2").Value = 0.4410279939397697
$ws.Range("D2").Value = 0.07930899002765557
$ws.Range("E2").Value = 0.4211715630783601
$ws.Range("G2").Value = 0.8285168191641503
$ws.Range("H2").Value = 0.772730941403637
$ws.Range("I2").Value = 0.5334048023530187
$ws.Range("N2").Value = 0.889068435693577
$ws.Range("B3").Value = 1.080042981006329
$ws.Range("C3").Value = 0.3843995217006295
$ws.Range("D3").Value = 0.07180443182238605
$ws.Range("E3").Value = 0.3670622442441811
$ws.Range("G3").Value = 0.7849318891588837
$ws.Range("H3").Value = 0.759117781820521
$ws.Range("I3").Value = 0.5259301036038053
$ws.Range("N3").Value = 0.9053844807486157
$ws.Range("B4").Value = 0.9923884694594562
$ws.Range("C4").Value = 0.3496615900773463
$ws.Range("D4").Value = 0.06723879186561987
$ws.Range("E4").Value = 0.3339756669310816
$ws.Range("G4").Value = 0.7589509797061851
$ws.Range("H4").Value = 0.7513724930354613
$ws.Range("I4").Value = 0.5218294836030992
$ws.Range("N4").Value = 0.91588475578409
$ws.Range("B5").Value = 0.9567273147872015
$ws.Range("C5").Value = 0.3355123973026934
$ws.Range("D5").Value = 0.06538867410129967
$ws.Range("E5").Value = 0.3205239994476585
$ws.Range("G5").Value = 0.7485568321405935
$ws.Range("H5").Value = 0.7483692020607009
$ws.Range("I5").Value = 0.5202802314421646
$ws.Range("N5").Value = 0.9202845641080107
$ws.Range("B6").Value = 0.950809304110237
$ws.Range("C6").Value = 0.3331633076853109
$ws.Range("D6").Value = 0.06508208712033081
$ws.Range("E6").Value = 0.3182921684433211
$ws.Range("G6").Value = 0.7468424798055366
$ws.Range("H6").Value = 0.7478797116694977
$ws.Range("I6").Value = 0.5200303013714134
$ws.Range("N6").Value = 0.9210224376159887
$ws.Range("B7").Value = 0.9919072950909822
$ws.Range("C7").Value = 0.3494707432811879
$ws.Range("D7").Value = 0.067213798598587
$ws.Range("E7").Value = 0.3337941298937039
$ws.Range("G7").Value = 0.7588100218982277
$ws.Range("H7").Value = 0.7513313718411325
$ws.Range("I7").Value = 0.5218080982737376
$ws.Range("N7").Value = 0.9159436042954566
$ws.Range("B8").Value = 1.173784594056656
$ws.Range("C8").Value = 0.4214950162170794
$ws.Range("D8").Value = 0.0767125438797791
$ws.Range("E8").Value = 0.4024844484551977
$ws.Range("G8").Value = 0.8133249578644666
$ws.Range("H8").Value = 0.7679090947945895
$ws.Range("I8").Value = 0.5307253624615669
$ws.Range("N8").Value = 0.8945939068138296
$ws.Range("B9").Value = 1.532460675891627
$ws.Range("C9").Value = 0.5630555793379699
$ws.Range("D9").Value = 0.09568398586480953
$ws.Range("E9").Value = 0.5384140519905571
$ws.Range("G9").Value = 0.9265611969658494
$ws.Range("H9").Value = 0.8053424280979016
$ws.Range("I9").Value = 0.5521448529413107
$ws.Range("N9").Value = 0.8565695002148395
$ws.Range("B10").Value = 1.797379534263314
$ws.Range("C10").Value = 0.667357868110571
$ws.Range("D10").Value = 0.10984720714562
$ws.Range("E10").Value = 0.6392386904091296
$ws.Range("G10").Value = 1.013820902331446
$ws.Range("H10").Value = 0.8359326481020446
$ws.Range("I10").Value = 0.5703565829551849
$ws.Range("N10").Value = 0.8309974018753611
$ws.Range("B11").Value = 1.918241302254387
$ws.Range("C11").Value = 0.714895335421545
$ws.Range("D11").Value = 0.1163426159284455
$ws.Range("E11").Value = 0.6853606222822037
$ws.Range("G11").Value = 1.05444556221758
$ws.Range("H11").Value = 0.8505385488238062
$ws.Range("I11").Value = 0.5791960274693508
$ws.Range("N11").Value = 0.8198822966329333
$ws.Range("B12").Value = 1.964061380303974
$ws.Range("C12").Value = 0.7329111484910982
$ws.Range("D12").Value = 0.1188100529928846
$ws.Range("E12").Value = 0.7028663188493312
$ws.Range("G12").Value = 1.069966323637743
$ws.Range("H12").Value = 0.856170188292765
$ws.Range("I12").Value = 0.5826244473903586
$ws.Range("N12").Value = 0.8157482068192727
$ws.Range("B13").Value = 1.954190842905291
$ws.Range("C13").Value = 0.7290304560642085
$ws.Range("D13").Value = 0.1182782973164649
$ws.Range("E13").Value = 0.6990943046543947
$ws.Range("G13").Value = 1.066617498697298
$ws.Range("H13").Value = 0.8549528119704064
$ws.Range("I13").Value = 0.5818824468724699
$ws.Range("N13").Value = 0.8166352150288926
$ws.Range("B14").Value = 1.922009885734894
$ws.Range("C14").Value = 0.7163772098502363
$ws.Range("D14").Value = 0.1165454565348512
$ws.Range("E14").Value = 0.6867999982667072
$ws.Range("G14").Value = 1.055719697670497
$ws.Range("H14").Value = 0.8509998406211139
$ws.Range("I14").Value = 0.5794764523107574
$ws.Range("N14").Value = 0.8195406773492004
$ws.Range("B15").Value = 1.902305003708307
$ws.Range("C15").Value = 0.7086286519760279
$ws.Range("D15").Value = 0.115485059546387
$ws.Range("E15").Value = 0.6792747357431921
$ws.Range("G15").Value = 1.049062433558817
$ws.Range("H15").Value = 0.8485916891503109
$ws.Range("I15").Value = 0.5780133143966353
$ws.Range("N15").Value = 0.8213301340002593
$ws.Range("B16").Value = 1.789488138552315
$ws.Range("C16").Value = 0.6642531169682115
$ws.Range("D16").Value = 0.1094237917330929
$ws.Range("E16").Value = 0.6362299659068356
$ws.Range("G16").Value = 1.011184982908389
$ws.Range("H16").Value = 0.8349921307303418
$ws.Range("I16").Value = 0.5697901824027838
$ws.Range("N16").Value = 0.8317342618715833
$ws.Range("B17").Value = 1.720369529037839
$ws.Range("C17").Value = 0.6370542206893219
$ws.Range("D17").Value = 0.1057189983310565
$ws.Range("E17").Value = 0.6098914838607499
$ws.Range("G17").Value = 0.9881889246686626
$ws.Range("H17").Value = 0.8268270529665074
$ws.Range("I17").Value = 0.5648886041664127
$ws.Range("N17").Value = 0.8382498260542472
$ws.Range("B18").Value = 1.68064703372454
$ws.Range("C18").Value = 0.6214184995207006
$ws.Range("D18").Value = 0.1035930262736144
$ws.Range("E18").Value = 0.5947661983524881
$ws.Range("G18").Value = 0.9750495104519246
$ws.Range("H18").Value = 0.8221955985921738
$ws.Range("I18").Value = 0.5621214893018376
$ws.Range("N18").Value = 0.8420461103836836
$ws.Range("B19").Value = 1.667203230221503
$ws.Range("C19").Value = 0.6161259048061538
$ws.Range("D19").Value = 0.1028740482975508
$ws.Range("E19").Value = 0.5896490390652076
$ws.Range("G19").Value = 0.9706156239115842
$ws.Range("H19").Value = 0.8206385677235062
$ws.Range("I19").Value = 0.5611935072141634
$ws.Range("N19").Value = 0.8433398168170955
$ws.Range("B20").Value = 1.727723922989185
$ws.Range("C20").Value = 0.6399487132497939
$ws.Range("D20").Value = 0.1061128684014108
$ws.Range("E20").Value = 0.6126927604561843
$ws.Range("G20").Value = 0.9906278337006142
$ws.Range("H20").Value = 0.8276895139718476
$ws.Range("I20").Value = 0.5654049793378562
$ws.Range("N20").Value = 0.8375511888153433
$ws.Range("B21").Value = 1.931460774437483
$ws.Range("C21").Value = 0.7200933719938689
$ws.Range("D21").Value = 0.1170542213262564
$ws.Range("E21").Value = 0.6904100117922525
$ws.Range("G21").Value = 1.058916898359826
$ws.Range("H21").Value = 0.8521581794785789
$ws.Range("I21").Value = 0.5801809391328732
$ws.Range("N21").Value = 0.8186852341817297
$ws.Range("B22").Value = 2.064920779466661
$ws.Range("C22").Value = 0.772557221735326
$ws.Range("D22").Value = 0.1242504028945461
$ws.Range("E22").Value = 0.7414394118619612
$ws.Range("G22").Value = 1.104348388317277
$ws.Range("H22").Value = 0.8687375219533351
$ws.Range("I22").Value = 0.5903113087821339
$ws.Range("N22").Value = 0.8067924049461229
$ws.Range("B23").Value = 1.993661906867487
$ws.Range("C23").Value = 0.7445480292867046
$ws.Range("D23").Value = 0.1204054401535757
$ws.Range("E23").Value = 0.7141812413952664
$ws.Range("G23").Value = 1.080026357828899
$ws.Range("H23").Value = 0.8598345609016462
$ws.Range("I23").Value = 0.5848607771707464
$ws.Range("N23").Value = 0.8130996615420404
$ws.Range("B24").Value = 1.724398955816412
$ws.Range("C24").Value = 0.6386401092678966
$ws.Range("D24").Value = 0.1059347874151371
$ws.Range("E24").Value = 0.6114262504347607
$ws.Range("G24").Value = 0.9895249498702299
$ws.Range("H24").Value = 0.8272993999842981
$ws.Range("I24").Value = 0.5651713677815309
$ws.Range("N24").Value = 0.8378668857654689
$ws.Range("B25").Value = 1.435196092539115
$ws.Range("C25").Value = 0.5247148793080783
$ws.Range("D25").Value = 0.09051330460006568
$ws.Range("E25").Value = 0.5014886775096841
$ws.Range("G25").Value = 0.8952288673592932
$ws.Range("H25").Value = 0.7946797264428653
$ws.Range("I25").Value = 0.5459213548925845
$ws.Range("N25").Value = 0.8664425269104403
